$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet updates (capital/pnl/trade-count roll-up after new trade #16
# closes on the MarketMaking strategy)
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.35   # Current Capital
$summary.Range("B4").Value = 0.35      # Total P&L $
$summary.Range("B5").Value = 0.44      # Total P&L %
$summary.Range("B6").Value = 16        # Total Trades
$summary.Range("B7").Value = 10        # Winning Trades
$summary.Range("B9").Value = 62.5      # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet updates (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.35     # Capital
$status.Range("D6").Value = 16         # Trades
$status.Range("E6").Value = 0.35       # P&L $
$status.Range("F6").Value = 0.35       # P&L %
$status.Range("G6").Value = 62.5       # Win Rate %

# ---------------------------------------------------------------------------
# Append the new closed trade (#16) as row 17 on both the "All Trades" log
# and the per-strategy "MarketMaking" log - the two sheets mirror each
# other row-for-row.
# ---------------------------------------------------------------------------
function Add-Trade17($ws) {
    $ws.Range("A17").Value = 16

    # "2026-02-17" parses as a real date through plain Value assignment, so
    # round it through a text formula first and flatten it back down to a
    # literal value/shared string via copy / paste-values - this keeps the
    # cell a plain string (matching the rest of the sheet) instead of
    # turning it into a styled date serial number.
    $ws.Range("B17").Formula = '="2026-02-17"'
    $ws.Range("B17").Copy()
    $ws.Range("B17").PasteSpecial(-4163)
    $excel.CutCopyMode = $false

    $ws.Range("C17").Value = "23:53:56"
    $ws.Range("D17").Value = "MarketMaking"
    $ws.Range("E17").Value = "UP"
    $ws.Range("F17").Value = 0.9399999999999999
    $ws.Range("G17").Value = 0.95
    $ws.Range("H17").Value = "CLOSED"
    $ws.Range("I17").Value = 1.0638
    $ws.Range("J17").Value = 0.01
    $ws.Range("K17").Value = 100.35
    $ws.Range("L17").Value = 0
    $ws.Range("M17").Value = 0
    $ws.Range("N17").Value = 0.6
    $ws.Range("O17").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P17").Value = "early_exit"
    $ws.Range("Q17").Value = 0.12
}

Add-Trade17 $wb.Worksheets.Item("All Trades")
Add-Trade17 $wb.Worksheets.Item("MarketMaking")
